# Add a new "Save" column (H) to the s_vals sheet, mirroring the header
# style used by the existing "sum" column (G) and filling the data rows
# with 0.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone G1's formatting (bold font, border, centered alignment) onto H1
# so the new header matches the look of the other header cells exactly.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# New data column values.
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
